# "Generate Report for Archive"
#
# The localization status for the two handed-off files moved on from
# "Ready for handoff" to "In Translation". That status string shows up
# in the "zh-cn" / "de-de" columns of the Overview sheet and in the
# "Status" column of each per-locale detail sheet, so update it
# everywhere it appears.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("Ready for handoff", "In Translation")
}

# The status column narrows to fit the new (shorter) text.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E:F").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C:C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C:C").ColumnWidth = 12.5
